{"js": "// Prepend \"Design: \" to each of the six \"List Bullet\" answer paragraphs\n// in the feedback table (the Question/Answers table under NB518:).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text,items/style\");\nawait context.sync();\n\nfor (let i = 0; i < paras.items.length; i++) {\n  const para = paras.items[i];\n  if (para.style === \"List Bullet\" && para.text && para.text.length > 0) {\n    para.getRange(\"Start\").insertText(\"Design: \", \"Start\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Prepend \"Design: \" to each of the six \"List Bullet\" answer paragraphs\n# in the feedback table (the Question/Answers table under NB518:).\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Range.Style.NameLocal\n    if ($styleName -eq \"List Bullet\" -and $p.Range.Text.Trim().Length -gt 0) {\n        $p.Range.InsertBefore(\"Design: \")\n    }\n}\n"}
